# Link from pregnancy to new child
# -------------------------------------------------------------------
# 1) survey sheet: duplicate the "if ... begin screen / note (linked_table
#    child) / end screen" block twice (LITTERSIZE >= '1' and >= '2'),
#    closing the newly nested ifs before the existing outer "end if".
# 2) queries sheet: repoint the existing SES linked_table query at
#    opendatakit.getCurrentInstanceId() instead of data('PREGID'), and add
#    a new "child" / CRIANCA linked_table query using the same id plus the
#    visit's outcome date.
# -------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# survey sheet (tab 2)
# ---------------------------------------------------------------
$survey = $wb.Worksheets.Item(2)

# Insert 12 fresh rows right before the old row 131 ("end if"), which
# pushes it down to row 143. The new rows inherit row 130's formatting
# (style 11 on columns A/B) the way Excel normally copies format from
# the row above on insert.
$survey.Range("A131:A142").EntireRow.Insert()
$survey.Range("A131:H142").ClearContents()

# give row 130 (now followed by the new "if" block) empty, styled G/H
# cells like the rest of the block below, matching the surrounding rows
$survey.Cells.Item(130, 7).Interior.Pattern = -4142
$survey.Cells.Item(130, 8).Interior.Pattern = -4142

# First nested if: data('LITTERSIZE') >='1'
$survey.Cells.Item(131, 2).Value = "if"
$survey.Cells.Item(131, 3).Value = "data('LITTERSIZE') >='1'"
$survey.Cells.Item(132, 2).Value = "begin screen"
$survey.Cells.Item(133, 4).Value = "note"
$survey.Cells.Item(134, 4).Value = "linked_table"
$survey.Cells.Item(134, 5).Value = "child"
$survey.Cells.Item(135, 2).Value = "end screen"

# Second nested if: data('LITTERSIZE') >='2'
$survey.Cells.Item(136, 2).Value = "if"
$survey.Cells.Item(136, 3).Value = "data('LITTERSIZE') >='2'"
$survey.Cells.Item(137, 2).Value = "begin screen"
$survey.Cells.Item(138, 4).Value = "note"
$survey.Cells.Item(139, 4).Value = "linked_table"
$survey.Cells.Item(139, 5).Value = "child"
$survey.Cells.Item(140, 2).Value = "end screen"

# close both new ifs (old row 131 "end if" - now row 143 - closes the
# outer if that was already there)
$survey.Cells.Item(141, 2).Value = "end if"
$survey.Cells.Item(142, 2).Value = "end if"

# empty styled G/H helper cells down the whole new block (rows 130-141)
foreach ($r in 130..141) {
    $survey.Cells.Item($r, 7).Interior.Pattern = -4142
    $survey.Cells.Item($r, 8).Interior.Pattern = -4142
}

# ---------------------------------------------------------------
# queries sheet (tab 3)
# ---------------------------------------------------------------
$queries = $wb.Worksheets.Item(3)

# the existing "ses" query now grabs the id of the pregnancy instance
# being saved, instead of re-reading it back out of its own field
$queries.Cells.Item(4, 6).Value = "[opendatakit.getCurrentInstanceId()]"
$queries.Cells.Item(4, 7).Value = "{PREGID: opendatakit.getCurrentInstanceId(), REGID: data('REGID'), SESDATA: data('PREGDIA')}"

# brand new query: link to the CRIANCA (child) table for this visit
$queries.Cells.Item(5, 1).Value = "child"
$queries.Cells.Item(5, 2).Value = "linked_table"
$queries.Cells.Item(5, 3).Value = "CRIANCA"
$queries.Cells.Item(5, 4).Value = "CRIANCA"
$queries.Cells.Item(5, 5).Value = "PREGID = ?"
$queries.Cells.Item(5, 6).Value = "[opendatakit.getCurrentInstanceId()]"
$queries.Cells.Item(5, 7).Value = "{PREGID: opendatakit.getCurrentInstanceId(), REGID: data('REGID'), REGDIA: data('PREGDIA'), OUTDATE: data('OUTDATE')}"
$queries.Cells.Item(5, 8).Value = "{}"

# ---------------------------------------------------------------
# view state: selections + which sheet/tab is active
# ---------------------------------------------------------------
$choices = $wb.Worksheets.Item(5)
$choices.Activate()
$choices.Range("A21:D23").Select()

$survey.Activate()
$survey.Range("F3").Select()

$queries.Activate()
$queries.Range("E8").Select()
